$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Select()
$excel.Goto($ws2.Range("M1"), $true)
